$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Data for columns D (label "U"), F and G, rows 2-11
$data = @(
    @(2, 1, 1),
    @(3, 3, 1),
    @(4, 5, 1),
    @(5, 3, 5),
    @(6, 5, 5),
    @(7, 6, 6),
    @(8, 8, 6),
    @(9, 10, 6),
    @(10, 8, 10),
    @(11, 10, 10)
)

foreach ($row in $data) {
    $r = $row[0]
    $fval = $row[1]
    $gval = $row[2]
    $ws.Cells.Item($r, 4).Value = "U"
    $ws.Cells.Item($r, 6).Value = $fval
    $ws.Cells.Item($r, 7).Value = $gval
}

# Update the selection / active cell range shown in sheetView
$ws.Range("D2:G11").Select()
